$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.579.53"
$ws.Range("E2").Value = "  +3.93%  "
$ws.Range("D3").Value = "1.742.63"
$ws.Range("E3").Value = "  +4.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.80%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4814"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2694"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("D10").Value = "1.743.73"
$ws.Range("E10").Value = "  +4.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07133"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.69%  "
$ws.Range("E14").Value = "  +3.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "26.585.91"
$ws.Range("E17").Value = "  +3.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006902"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("D21").Value = "1.968.31"
$ws.Range("E21").Value = "  +4.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.628"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.843"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.364"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.818"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.433"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.004"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.739"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07887"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04609"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.619"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6388"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9964"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9363"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "113.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.000"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.435"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.003"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.766"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +18.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01507"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3915"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1226"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.746"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05337"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.959"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.262"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3452"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.24%  "
